$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.769.93'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '3.530.38'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''626.16'
$ws.Range('E5').Value = '  +3.14%  '
$ws.Range('D6').Value = '''174.30'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').Value = '3.526.20'
$ws.Range('E7').Value = '  -1.00%  '
$ws.Range('D8').Value = '''0.609'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').Value = '''1.00'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '''0.197'
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('E11').Value = '  -2.99%  '
$ws.Range('D12').Value = '''0.586'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '''46.50'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').Value = '''0.0000276'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '4.098.91'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').Value = '''8.41'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '''607.04'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('D18').Value = '3.527.71'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').Value = '70.866.97'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('D21').Value = '''17.79'
$ws.Range('E21').Value = '  +2.09%  '
$ws.Range('D22').Value = '''0.885'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '''9.08'
$ws.Range('E23').Value = '  -3.62%  '
$ws.Range('D24').Value = '''15.66'
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('D25').Value = '''98.15'
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '''2.58'
$ws.Range('E28').Value = '  -2.27%  '
$ws.Range('D29').Value = '''33.81'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('D30').Value = '''9.09'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = '''8.13'
$ws.Range('E32').Value = '  -4.15%  '
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = '''640.03'
$ws.Range('E34').Value = '  +4.66%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '''6.82'
$ws.Range('E35').Value = '  -2.84%  '
$ws.Range('D36').Value = '''0.0999'
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('D37').Value = '''10.83'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').Value = '''3.47'
$ws.Range('E38').Value = '  -6.65%  '
$ws.Range('D39').Value = '''0.0475'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').Value = '''56.79'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').Value = '''0.143'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('D43').Value = '3.359.55'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').Value = '0.0₃0734'
$ws.Range('E44').Value = '  +3.71%  '
$ws.Range('D45').Value = '''3.00'
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('D46').Value = '''0.312'
$ws.Range('E46').Value = '  -2.74%  '
$ws.Range('D47').Value = '''32.10'
$ws.Range('E47').Value = '  -2.89%  '
$ws.Range('D48').Value = '''2.55'
$ws.Range('E48').Value = '  -2.06%  '
$ws.Range('D49').Value = '''0.130'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '''132.89'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('E51').Value = '  +5.90%  '
